$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the oldest data point), shifting all subsequent rows up by one.
# This naturally drops the old last row (19) off the bottom, matching the new
# A1:E18 dimension, and carries column A's date style upward with the cells.
$ws.Rows.Item(2).Delete() | Out-Null

# Recalculated y_1_forecast (column E) values for the now-shifted rows 2:18.
$newE = @(
  2.651903832953884,
  0.5356365903280924,
  1.805987045940682,
  4.382531137514767,
  1.988448192515935,
  3.166945525867848,
  3.684750195712683,
  2.9587404276884,
  2.661643377950096,
  2.433248629349549,
  3.296423324101938,
  2.763966172716947,
  1.856930494010856,
  4.270817433327112,
  3.604316462518464,
  2.098908173995873,
  1.874466487556892
)

for ($i = 0; $i -lt $newE.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $newE[$i]
}
